$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Individual cell value edits that do NOT involve row shifting (rows 2-25) ---
# D5 ("RM 14"): -14.4 -> blank
$ws.Range("D5").Value = ""

# D11 ("RM 58"): blank -> -15.5
$ws.Range("D11").Value = -15.5

# C19 ("RM 125"): blank -> 13.2
$ws.Range("C19").Value = 13.2
# D19 ("RM 125"): -15.5 -> blank
$ws.Range("D19").Value = ""

# C21 ("RM 135"): 12.7 -> blank
$ws.Range("C21").Value = ""

# C23 ("RM 140"): blank -> 12.2
$ws.Range("C23").Value = 12.2

# D25 ("RM 145"): blank -> -15.5
$ws.Range("D25").Value = -15.5

# --- Delete entire rows (original numbering). Delete higher row index first ---
# Row 28 = "SC 92" is removed entirely
$ws.Rows.Item(28).Delete()
# Row 26 = "RM 232" is removed entirely
$ws.Rows.Item(26).Delete()

# --- After the two row deletions, rows below shifted up by the amounts noted.
# Fix remaining individual value diffs using the NEW (post-deletion) row numbers ---
# New row 27 ("SC 101"): C was 10 -> blank
$ws.Range("C27").Value = ""

# New row 29 ("SC 119"): D was -13 -> blank
$ws.Range("D29").Value = ""

# New row 33 ("SC 232"): C was blank -> 10.4
$ws.Range("C33").Value = 10.4
